$wb = $excel.ActiveWorkbook

# Sheet "建物" (building) - property_category column I, rows 2-5: land -> building
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I5").Value = "building"

# Sheet "汽車" (car) - property_category column H, rows 2-3: land -> car
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2:H3").Value = "car"
